$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K = "Дата передачи показаний" (Reading submission date).
# Rows 5, 9 and 10 get their date value changed to the short "YYYY-MM" form "2021-12"
# to reflect that the parser now also accepts dates given as just year-month.
$ws.Range("K5").Value = "2021-12"
$ws.Range("K9").Value = "2021-12"
$ws.Range("K10").Value = "2021-12"
